# Generate Report for Handoff
# Update the localization-status report: both the zh-cn and de-de locales
# moved from "In Translation" to "Ready for handoff", their handoff-xliff
# generation timestamps were refreshed, and the "Status" columns were
# widened to fit the new (longer) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"       # Status column
$wsDeDe.Range("C2").Value = "Ready for handoff"       # Status column

# --- Refreshed handoff timestamps ------------------------------------------
# de-de handoff datetime (shared between the Overview sheet and the de-de
# sheet, same as in the original workbook).
$wsOverview.Range("G2").Value = "2016-09-04 03:03:09"
$wsDeDe.Range("H2").Value = "2016-09-04 03:03:09"

# zh-cn handoff datetime.
$wsZhCn.Range("H2").Value = "2016-09-04 03:02:59"

# --- Widen the Status columns to fit the new text ---------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 16.3   # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3        # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3        # column C (Status)
